$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy style from existing header cell (H1) to new header cells so they match (bold, centered, bordered)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Data values for columns I (I0) and J (IF)
$data = @{
    2  = @(10, 10)
    3  = @(6, 6)
    4  = @(5, 6)
    5  = @(6, 6)
    6  = @(9, 9)
    7  = @(8, 8)
    8  = @(6, 6)
    9  = @(7, 7)
    10 = @(8, 8)
    11 = @(7, 7)
    12 = @(7, 7)
    13 = @(6, 7)
    14 = @(8, 8)
    15 = @(8, 9)
    16 = @(9, 9)
    17 = @(9, 9)
    18 = @(6, 6)
    19 = @(9, 9)
    20 = @(8, 9)
    21 = @(9, 9)
    22 = @(5, 5)
    23 = @(4, 4)
    24 = @(5, 5)
    25 = @(9, 9)
    26 = @(6, 6)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
